# Apply the "Updated cryptos list" data refresh to Sheet1.
# Column D holds price text (kept as text via a leading apostrophe so
# Excel does not reinterpret values such as "324.67" or "28.903.93" as numbers),
# column E holds the padded percentage text, and a few rows (44/45, 47/48)
# had their coin data reordered.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range('D2').Value = '''28.903.93'
$ws.Range('E2').Value = '  -1.27%  '

# Row 3 (Ethereum)
$ws.Range('D3').Value = '''1.908.29'
$ws.Range('E3').Value = '  -1.23%  '

# Row 4 (TetherUSD)
$ws.Range('E4').Value = '  +0.25%  '

# Row 5 (BNB)
$ws.Range('D5').Value = '''324.67'
$ws.Range('E5').Value = '  -0.21%  '

# Row 6 (USDC)
$ws.Range('E6').Value = '  +0.12%  '

# Row 7 (XRP)
$ws.Range('D7').Value = '''0.4582'
$ws.Range('E7').Value = '  -0.85%  '

# Row 8 (Cardano)
$ws.Range('D8').Value = '''0.3815'
$ws.Range('E8').Value = '  -1.39%  '

# Row 9 (Dogecoin)
$ws.Range('D9').Value = '''0.07729'
$ws.Range('E9').Value = '  -1.06%  '

# Row 10 (Polygon)
$ws.Range('D10').Value = '''0.9796'
$ws.Range('E10').Value = '  +0.82%  '

# Row 11 (Solana)
$ws.Range('D11').Value = '''22.08'
$ws.Range('E11').Value = '  -2.23%  '

# Row 12 (WrappedEther)
$ws.Range('D12').Value = '''1.897.59'
$ws.Range('E12').Value = '  -2.21%  '

# Row 13 (Polkadot)
$ws.Range('E13').Value = '  -1.60%  '

# Row 14 (Chainlink)
$ws.Range('D14').Value = '''6.942'
$ws.Range('E14').Value = '  -1.83%  '

# Row 15 (TRON)
$ws.Range('D15').Value = '''0.07067'
$ws.Range('E15').Value = '  +0.06%  '

# Row 16 (BinanceUSD)
$ws.Range('E16').Value = '  +0.11%  '

# Row 17 (Litecoin)
$ws.Range('D17').Value = '''83.71'
$ws.Range('E17').Value = '  -3.46%  '

# Row 18 (ShibaInu)
$ws.Range('D18').Value = '''0.000009458'
$ws.Range('E18').Value = '  -2.22%  '

# Row 19 (Avalanche)
$ws.Range('E19').Value = '  -2.42%  '

# Row 20 (Dai)
$ws.Range('E20').Value = '  +0.00%  '

# Row 21 (WrappedBTC)
$ws.Range('D21').Value = '''28.897.57'
$ws.Range('E21').Value = '  -1.40%  '

# Row 22 (Uniswap)
$ws.Range('D22').Value = '''5.315'
$ws.Range('E22').Value = '  -2.86%  '

# Row 23 (Cosmos)
$ws.Range('D23').Value = '''10.93'
$ws.Range('E23').Value = '  -1.05%  '

# Row 24 (Toncoin)
$ws.Range('D24').Value = '''2.098'
$ws.Range('E24').Value = '  +0.39%  '

# Row 25 (Monero)
$ws.Range('D25').Value = '''158.63'
$ws.Range('E25').Value = '  +1.05%  '

# Row 26 (EthereumClassic)
$ws.Range('D26').Value = '''19.07'
$ws.Range('E26').Value = '  -1.40%  '

# Row 27 (InternetComputer(DFINITY))
$ws.Range('E27').Value = '  -1.48%  '

# Row 28 (BitcoinCash)
$ws.Range('D28').Value = '''117.45'
$ws.Range('E28').Value = '  -0.71%  '

# Row 29 (LidoDAOToken)
$ws.Range('D29').Value = '''1.870'
$ws.Range('E29').Value = '  +0.75%  '

# Row 30 (Stellar)
$ws.Range('E30').Value = '  -0.43%  '

# Row 31 (ImmutableX)
$ws.Range('E31').Value = '  -0.13%  '

# Row 32 (Filecoin)
$ws.Range('D32').Value = '''5.087'
$ws.Range('E32').Value = '  -1.48%  '

# Row 33 (ARBITRUM)
$ws.Range('E33').Value = '  -4.42%  '

# Row 34 (HuobiToken)
$ws.Range('D34').Value = '''3.037'
$ws.Range('E34').Value = '  -1.18%  '

# Row 35 (Hedera)
$ws.Range('D35').Value = '''0.05708'
$ws.Range('E35').Value = '  -1.08%  '

# Row 36 (TrustWalletToken)
$ws.Range('D36').Value = '''1.156'
$ws.Range('E36').Value = '  +0.13%  '

# Row 37 (Frax)
$ws.Range('D37').Value = '''1.001'
$ws.Range('E37').Value = '  +0.20%  '

# Row 38 (VeChain)
$ws.Range('D38').Value = '''0.02043'
$ws.Range('E38').Value = '  -1.92%  '

# Row 39 (FraxShare)
$ws.Range('D39').Value = '''7.418'

# Row 40 (TheSandbox)
$ws.Range('E40').Value = '  -2.91%  '

# Row 41 (Algorand)
$ws.Range('E41').Value = '  -1.45%  '

# Row 42 (MXToken)
$ws.Range('D42').Value = '''2.859'
$ws.Range('E42').Value = '  +5.13%  '

# Row 43 (Aptos)
$ws.Range('D43').Value = '''9.313'
$ws.Range('E43').Value = '  -0.53%  '

# Row 44 (RenderToken)
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '''0.5172'
$ws.Range('E44').Value = '  -1.87%  '

# Row 45 (Decentraland)
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''2.126'
$ws.Range('E45').Value = '  +2.44%  '

# Row 46 (EnergySwap)
$ws.Range('D46').Value = '''11.21'
$ws.Range('E46').Value = '  -2.73%  '

# Row 47 (Cronos)
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '''0.000002661'
$ws.Range('E47').Value = '  -13.96%  '

# Row 48 (PEPE)
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.06901'
$ws.Range('E48').Value = '  +0.56%  '

# Row 49 (NEARProtocol)
$ws.Range('D49').Value = '''1.777'
$ws.Range('E49').Value = '  -1.81%  '

# Row 50 (Quant)
$ws.Range('D50').Value = '''110.32'
$ws.Range('E50').Value = '  -0.86%  '

# Row 51 (WOONetwork)
$ws.Range('D51').Value = '''0.2882'
$ws.Range('E51').Value = '  -3.71%  '
